$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 248.7642878020601

$ws.Range("B3").Value = 0.006876353814593728
$ws.Range("C3").Value = 86.29678392075563
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 71517.89157740913
$ws.Range("G3").Value = 71604.90797053471
